$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list data (price + 1h volume change columns, and a few
# coin name/link swaps) to match the latest scrape, cell by cell.

$ws.Range("D2").Value = "'48.666.53"
$ws.Range("E2").Value = "'  +8.35%  "
$ws.Range("D3").Value = "'2.668.30"
$ws.Range("E3").Value = "'  +11.91%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'314.88"
$ws.Range("E5").Value = "'  +7.74%  "
$ws.Range("D6").Value = "'106.45"
$ws.Range("E6").Value = "'  +14.91%  "
$ws.Range("D7").Value = "'0.616"
$ws.Range("E7").Value = "'  +11.09%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("D9").Value = "'0.602"
$ws.Range("E9").Value = "'  +21.47%  "
$ws.Range("D10").Value = "'41.59"
$ws.Range("E10").Value = "'  +22.86%  "
$ws.Range("D11").Value = "'0.0876"
$ws.Range("E11").Value = "'  +13.55%  "
$ws.Range("D12").Value = "'55.87"
$ws.Range("E12").Value = "'  +4.74%  "
$ws.Range("D13").Value = "'8.49"
$ws.Range("E13").Value = "'  +22.20%  "
$ws.Range("D14").Value = "'3.081.63"
$ws.Range("E14").Value = "'  +11.98%  "
$ws.Range("E15").Value = "'  +3.91%  "
$ws.Range("D16").Value = "'2.658.13"
$ws.Range("E16").Value = "'  +11.24%  "
$ws.Range("D17").Value = "'0.952"
$ws.Range("E17").Value = "'  +16.02%  "
$ws.Range("D18").Value = "'15.54"
$ws.Range("E18").Value = "'  +11.09%  "
$ws.Range("D19").Value = "'48.803.04"
$ws.Range("E19").Value = "'  +8.55%  "
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("E20").Value = "'  +12.62%  "
$ws.Range("D21").Value = "'13.46"
$ws.Range("E21").Value = "'  +9.43%  "
$ws.Range("D22").Value = "'6.93"
$ws.Range("E22").Value = "'  +14.42%  "
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'285.86"
$ws.Range("E23").Value = "'  +20.67%  "
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'73.95"
$ws.Range("E24").Value = "'  +11.52%  "
$ws.Range("E25").Value = "'  +13.83%  "
$ws.Range("E26").Value = "'  +20.17%  "
$ws.Range("D27").Value = "'30.50"
$ws.Range("E27").Value = "'  +46.92%  "
$ws.Range("E28").Value = "'  +0.24%  "
$ws.Range("E29").Value = "'  +2.80%  "
$ws.Range("D30").Value = "'10.87"
$ws.Range("E30").Value = "'  +14.68%  "
$ws.Range("D31").Value = "'41.18"
$ws.Range("E31").Value = "'  +10.25%  "
$ws.Range("E32").Value = "'  +5.05%  "
$ws.Range("D33").Value = "'6.31"
$ws.Range("E33").Value = "'  +18.07%  "
$ws.Range("D34").Value = "'3.74"
$ws.Range("E34").Value = "'  -0.96%  "
$ws.Range("D35").Value = "'0.0871"
$ws.Range("E35").Value = "'  +15.79%  "
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'2.27"
$ws.Range("E36").Value = "'  +16.51%  "
$ws.Range("B37").Value = "'WEMIXToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'2.89"
$ws.Range("E37").Value = "'  +7.01%  "
$ws.Range("B38").Value = "'Kaspa"
$ws.Range("C38").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.128"
$ws.Range("E38").Value = "'  +15.02%  "
$ws.Range("B39").Value = "'Monero"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'153.87"
$ws.Range("E39").Value = "'  +4.57%  "
$ws.Range("E40").Value = "'  +10.37%  "
$ws.Range("D41").Value = "'16.67"
$ws.Range("E41").Value = "'  +16.20%  "
$ws.Range("D42").Value = "'4.38"
$ws.Range("E42").Value = "'  +19.04%  "
$ws.Range("B43").Value = "'NEARProtocol"
$ws.Range("C43").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.81"
$ws.Range("E43").Value = "'  +22.08%  "
$ws.Range("B44").Value = "'EnergySwap"
$ws.Range("C44").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'23.10"
$ws.Range("E44").Value = "'  +50.40%  "
$ws.Range("D45").Value = "'0.0341"
$ws.Range("E45").Value = "'  +17.46%  "
$ws.Range("D46").Value = "'2.219.70"
$ws.Range("E46").Value = "'  +12.70%  "
$ws.Range("D47").Value = "'98.84"
$ws.Range("E47").Value = "'  +12.66%  "
$ws.Range("D48").Value = "'0.995"
$ws.Range("E48").Value = "'  -0.37%  "
$ws.Range("D49").Value = "'10.02"
$ws.Range("E49").Value = "'  +19.93%  "
$ws.Range("D50").Value = "'1.95"
$ws.Range("E50").Value = "'  +15.63%  "
$ws.Range("D51").Value = "'115.91"
$ws.Range("E51").Value = "'  +17.36%  "
